# penambahan field di mst022
# Adds a new "BLNC010 / Saldo Hutang Cancel Booking" row to the "saldo"
# sheet and makes that sheet the active tab (it was previously "Master").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("saldo")

# Fill column B (Kode Tabel) before column A (Nama Tabel) so the new
# shared-string entries are appended in the same order as the target file
# (BLNC010 first, then the longer description).
$ws.Range("B5").Value = "BLNC010"
$ws.Range("A5").Value = "Saldo Hutang Cancel Booking"

# Make "saldo" the active/selected sheet and move the selection to A6,
# just below the newly added row.
$ws.Activate()
$ws.Range("A6").Select()
